$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price/volume figures (GitHub Actions scheduled refresh).
# Price cells (column D) that are purely numeric-looking strings need to be
# forced back to Text so they keep matching the sheet's original inline-string
# cells (e.g. "505.91") instead of being auto-converted to numbers by Excel.

$ws.Range("D2").Value = "55.231.96"
$ws.Range("E2").Value = "  +1.31%  "
$ws.Range("D3").Value = "2.292.01"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "505.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.65"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.91%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.532"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.25%  "
$ws.Range("D9").Value = "2.312.83"
$ws.Range("E9").Value = "  +0.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0981"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.47%  "
$ws.Range("E11").Value = "  +1.73%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.09"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +7.30%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.340"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.77"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.78%  "
$ws.Range("D15").Value = "2.701.70"
$ws.Range("E15").Value = "  +0.06%  "
$ws.Range("D16").Value = "55.070.05"
$ws.Range("E16").Value = "  +1.07%  "
$ws.Range("E17").Value = "  +1.45%  "
$ws.Range("D18").Value = "2.294.14"
$ws.Range("E18").Value = "  -0.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.51"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.06%  "
$ws.Range("E20").Value = "  +0.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "310.61"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.75%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.59"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.71%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.50"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.34%  "
$ws.Range("E25").Value = "  -0.20%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("E27").Value = "  +2.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "173.12"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.05%  "
$ws.Range("E29").Value = "  +3.04%  "
$ws.Range("D30").Value = "0.0₃0709"
$ws.Range("E30").Value = "  +1.41%  "
$ws.Range("E31").Value = "  +0.28%  "
$ws.Range("E32").Value = "  +4.35%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E34").Value = "  +1.12%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.997"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.922"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.49%  "
$ws.Range("E37").Value = "  +2.35%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.88"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.82"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.24%  "
$ws.Range("E40").Value = "  +2.05%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "133.56"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.64%  "
$ws.Range("E43").Value = "  +1.17%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.92"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "260.88"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.50%  "
$ws.Range("E46").Value = "  +1.96%  "
$ws.Range("E47").Value = "  +1.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.551"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.23%  "
$ws.Range("E49").Value = "  +0.35%  "
$ws.Range("E50").Value = "  +2.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.46"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.19%  "
